$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label swap for BP1/BQ1 (avoid shared-string collision by
#     writing the non-colliding text first) ---
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("BP1").Value = "average_doctor_old"

# --- Row-by-row data updates for the "average_doctor" columns (BP/BQ). ---
# The previous BP (new) values become the BQ (old) values, and BP gets
# freshly computed values.
$ws.Range("BQ4").Value = $ws.Range("BP4").Value2
$ws.Range("BP4").Value = 0.672
$ws.Range("BQ5").Value = $ws.Range("BP5").Value2
$ws.Range("BP5").Value = 0.449
$ws.Range("BQ6").Value = $ws.Range("BP6").Value2
$ws.Range("BP6").Value = 0.535
$ws.Range("BQ7").Value = $ws.Range("BP7").Value2
$ws.Range("BP7").Value = 0.479
$ws.Range("BQ8").Value = $ws.Range("BP8").Value2
$ws.Range("BP8").Value = 0.584
$ws.Range("BQ9").Value = $ws.Range("BP9").Value2
$ws.Range("BP9").Value = 0.585
$ws.Range("BQ10").Value = $ws.Range("BP10").Value2
$ws.Range("BP10").Value = 0.667
$ws.Range("BQ11").Value = $ws.Range("BP11").Value2
$ws.Range("BP11").Value = 0.667
$ws.Range("BQ12").Value = $ws.Range("BP12").Value2
$ws.Range("BP12").Value = 1.125
$ws.Range("BQ13").Value = $ws.Range("BP13").Value2
$ws.Range("BP13").Value = 0.795

# --- Plain numeric cell value updates for all other changed cells ---
$ws.Range("E4").Value = 0.393
$ws.Range("E5").Value = 0.493
$ws.Range("E6").Value = 0.437
$ws.Range("E7").Value = 0.469
$ws.Range("E8").Value = 0.51
$ws.Range("E9").Value = 0.422
$ws.Range("E10").Value = 0.556
$ws.Range("E11").Value = 0.556
$ws.Range("E12").Value = 1.36
$ws.Range("E13").Value = 1.676
$ws.Range("F4").Value = 0.089
$ws.Range("F5").Value = 0.104
$ws.Range("F8").Value = 0.134
$ws.Range("F9").Value = 0.244
$ws.Range("F10").Value = 0.247
$ws.Range("F11").Value = 0.247
$ws.Range("F12").Value = 0.47
$ws.Range("F13").Value = 0.743
$ws.Range("G4").Value = 0.299
$ws.Range("G5").Value = 0.322
$ws.Range("G8").Value = 0.367
$ws.Range("G9").Value = 0.494
$ws.Range("G10").Value = 0.497
$ws.Range("G11").Value = 0.497
$ws.Range("G12").Value = 0.6860000000000001
$ws.Range("G13").Value = 0.862
$ws.Range("N4").Value = 0.404
$ws.Range("N5").Value = 0.743
$ws.Range("N6").Value = 0.523
$ws.Range("N7").Value = 0.636
$ws.Range("N8").Value = 0.732
$ws.Range("N9").Value = 0.6
$ws.Range("N10").Value = 0.8
$ws.Range("N11").Value = 0.844
$ws.Range("N12").Value = 1.725
$ws.Range("N13").Value = 2.32
$ws.Range("O4").Value = 0.07000000000000001
$ws.Range("O5").Value = 0.09
$ws.Range("O8").Value = 0.077
$ws.Range("O9").Value = 0.24
$ws.Range("O10").Value = 0.16
$ws.Range("O11").Value = 0.131
$ws.Range("O12").Value = 1.799
$ws.Range("O13").Value = 1.259
$ws.Range("P4").Value = 0.264
$ws.Range("P5").Value = 0.3
$ws.Range("P8").Value = 0.278
$ws.Range("P9").Value = 0.49
$ws.Range("P10").Value = 0.4
$ws.Range("P11").Value = 0.362
$ws.Range("P12").Value = 1.341
$ws.Range("P13").Value = 1.122
$ws.Range("W4").Value = 0.215
$ws.Range("W5").Value = 0.217
$ws.Range("W6").Value = 0.216
$ws.Range("W7").Value = 0.217
$ws.Range("W8").Value = 0.216
$ws.Range("W9").Value = 0.111
$ws.Range("W10").Value = 0.267
$ws.Range("W11").Value = 0.267
$ws.Range("W12").Value = 1.917
$ws.Range("W13").Value = 1.073
$ws.Range("X4").Value = 0.103
$ws.Range("X5").Value = 0.108
$ws.Range("X8").Value = 0.1
$ws.Range("X9").Value = 0.099
$ws.Range("X10").Value = 0.196
$ws.Range("X11").Value = 0.196
$ws.Range("X12").Value = 0.743
$ws.Range("X13").Value = 0.173
$ws.Range("Y4").Value = 0.321
$ws.Range("Y5").Value = 0.328
$ws.Range("Y8").Value = 0.317
$ws.Range("Y9").Value = 0.314
$ws.Range("Y10").Value = 0.442
$ws.Range("Y11").Value = 0.442
$ws.Range("Y12").Value = 0.862
$ws.Range("Y13").Value = 0.416
$ws.Range("AI4").Value = 0.177
$ws.Range("AI5").Value = 0.217
$ws.Range("AI6").Value = 0.195
$ws.Range("AI7").Value = 0.208
$ws.Range("AI8").Value = 0.184
$ws.Range("AI9").Value = 0.089
$ws.Range("AI10").Value = 0.2
$ws.Range("AI11").Value = 0.2
$ws.Range("AI12").Value = 2.111
$ws.Range("AI13").Value = 1.385
$ws.Range("AJ4").Value = 0.05
$ws.Range("AJ5").Value = 0.083
$ws.Range("AJ8").Value = 0.073
$ws.Range("AJ9").Value = 0.081
$ws.Range("AJ10").Value = 0.16
$ws.Range("AJ11").Value = 0.16
$ws.Range("AJ12").Value = 0.988
$ws.Range("AJ13").Value = 0.435
$ws.Range("AK4").Value = 0.224
$ws.Range("AK5").Value = 0.289
$ws.Range("AK8").Value = 0.271
$ws.Range("AK9").Value = 0.285
$ws.Range("AK10").Value = 0.4
$ws.Range("AK11").Value = 0.4
$ws.Range("AK12").Value = 0.994
$ws.Range("AK13").Value = 0.66
$ws.Range("AU4").Value = 0.145
$ws.Range("AU5").Value = 0.29
$ws.Range("AU6").Value = 0.193
$ws.Range("AU7").Value = 0.242
$ws.Range("AU8").Value = 0.216
$ws.Range("AU10").Value = 0.178
$ws.Range("AU11").Value = 0.311
$ws.Range("AU12").Value = 2.929
$ws.Range("AU13").Value = 2.472
$ws.Range("AV4").Value = 0.028
$ws.Range("AV5").Value = 0.094
$ws.Range("AV8").Value = 0.067
$ws.Range("AV10").Value = 0.146
$ws.Range("AV11").Value = 0.214
$ws.Range("AV12").Value = 1.638
$ws.Range("AV13").Value = 1.397
$ws.Range("AW4").Value = 0.167
$ws.Range("AW5").Value = 0.307
$ws.Range("AW8").Value = 0.258
$ws.Range("AW10").Value = 0.382
$ws.Range("AW11").Value = 0.463
$ws.Range("AW12").Value = 1.28
$ws.Range("AW13").Value = 1.182
$ws.Range("BA4").Value = 2.017
$ws.Range("BA5").Value = 1.348
$ws.Range("BA6").Value = 1.604
$ws.Range("BA7").Value = 1.438
$ws.Range("BA8").Value = 1.753
$ws.Range("BA9").Value = 1.756
$ws.Range("BA10").Value = 2
$ws.Range("BA11").Value = 2
$ws.Range("BA12").Value = 3.375
$ws.Range("BA13").Value = 2.384
$ws.Range("BB4").Value = 0.174
$ws.Range("BB5").Value = 0.07000000000000001
$ws.Range("BB8").Value = 0.137
$ws.Range("BB9").Value = 0.249
$ws.Range("BB10").Value = 0.25
$ws.Range("BB11").Value = 0.25
$ws.Range("BB12").Value = 0.043
$ws.Range("BB13").Value = 0.292
$ws.Range("BC4").Value = 0.417
$ws.Range("BC5").Value = 0.265
$ws.Range("BC8").Value = 0.369
$ws.Range("BC9").Value = 0.499
$ws.Range("BC10").Value = 0.5
$ws.Range("BC11").Value = 0.5
$ws.Range("BC12").Value = 0.208
$ws.Range("BC13").Value = 0.54
$ws.Range("BG4").Value = 0.756
$ws.Range("BG5").Value = 0.407
$ws.Range("BG6").Value = 0.529
$ws.Range("BG7").Value = 0.448
$ws.Range("BG8").Value = 0.58
$ws.Range("BG9").Value = 0.622
$ws.Range("BG10").Value = 0.667
$ws.Range("BG11").Value = 0.667
$ws.Range("BG12").Value = 1.067
$ws.Range("BG13").Value = 0.579
$ws.Range("BH4").Value = 0.14
$ws.Range("BH5").Value = 0.046
$ws.Range("BH8").Value = 0.104
$ws.Range("BH9").Value = 0.235
$ws.Range("BH10").Value = 0.222
$ws.Range("BH11").Value = 0.222
$ws.Range("BH12").Value = 0.062
$ws.Range("BH13").Value = 0.044
$ws.Range("BI4").Value = 0.374
$ws.Range("BI5").Value = 0.215
$ws.Range("BI8").Value = 0.322
$ws.Range("BI9").Value = 0.485
$ws.Range("BI10").Value = 0.471
$ws.Range("BI11").Value = 0.471
$ws.Range("BI12").Value = 0.249
$ws.Range("BI13").Value = 0.21
$ws.Range("BM4").Value = 0.698
$ws.Range("BM5").Value = 0.571
$ws.Range("BM6").Value = 0.628
$ws.Range("BM7").Value = 0.593
$ws.Range("BM8").Value = 0.712
$ws.Range("BM9").Value = 0.667
$ws.Range("BM10").Value = 0.844
$ws.Range("BM11").Value = 0.844
$ws.Range("BM12").Value = 1.263
$ws.Range("BM13").Value = 0.958
$ws.Range("BN4").Value = 0.092
$ws.Range("BN5").Value = 0.081
$ws.Range("BN8").Value = 0.074
$ws.Range("BN9").Value = 0.222
$ws.Range("BN10").Value = 0.131
$ws.Range("BN11").Value = 0.131
$ws.Range("BN12").Value = 0.299
$ws.Range("BN13").Value = 0.302
$ws.Range("BO4").Value = 0.303
$ws.Range("BO5").Value = 0.284
$ws.Range("BO8").Value = 0.272
$ws.Range("BO9").Value = 0.471
$ws.Range("BO10").Value = 0.362
$ws.Range("BO11").Value = 0.362
$ws.Range("BO12").Value = 0.547
$ws.Range("BO13").Value = 0.55
